$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E27").Value = "[b'we']"
$ws.Range("F27").Value = 1

$ws.Range("E28").Value = "[b'scrutinize']"
$ws.Range("F28").Value = 1

$ws.Range("E29").Value = "[b'what']"
$ws.Range("F29").Value = 1

$ws.Range("F29").Select()
